$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the "Meta description: ..." paragraph that currently sits right
#    after the H1 title ("Play Champion Raceway Free: Review & Gameplay
#    Mechanics").
# ---------------------------------------------------------------------------
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "Meta description*") {
        $para.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# 2. Insert a new bold paragraph containing the title text right before the
#    final paragraph (the one that used to hold the image-prompt text), and
#    make sure the paragraph that precedes it ("Standard symbol payouts are
#    low") keeps its original list-bullet formatting.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$precedingPara = $d.Paragraphs.Item($count - 1)
$precedingText = $precedingPara.Range.Text.TrimEnd([char]13, [char]7)

$fragment = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
              '<w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr>' +
              '<w:r/>' +
              '<w:r><w:t>' + $precedingText + '</w:t></w:r>' +
            '</w:p>' +
            '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
              '<w:r/>' +
              '<w:r><w:rPr><w:b/></w:rPr><w:t>Play Champion Raceway Free: Review &amp; Gameplay Mechanics</w:t></w:r>' +
            '</w:p>'

$precedingPara.Range.InsertXML($fragment) | Out-Null

# ---------------------------------------------------------------------------
# 3. Swap the final paragraph's italic text (the old image-generation prompt)
#    for the meta-description copy.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Create a cartoon-style feature image for Champion Raceway that features a happy Maya warrior with glasses. The image should showcase the horse racing theme of the game, with the Maya warrior standing next to a horse on a racetrack. Use bright colors to make the image eye-catching and visually appealing. Add in elements from the game, such as the Wild symbol or the racetrack above the reels, to tie it back to the game. Overall, the feature image should capture the fun and excitement of horse racing while also highlighting the unique aspects of Champion Raceway.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Discover the gameplay mechanics, pros, and cons of Champion Raceway in our review. Play for free and experience the unique multiplier system and two free spin features.",
    2
) | Out-Null
